$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1160.8108
$ws.Range("J17").Value = 1236.7646
$ws.Range("L17").Value = 3710.2938
$ws.Range("N17").Value = -4046.2938
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984
$ws.Range("H132").Value = 22655556
$ws.Range("I132").Value = 31716428
$ws.Range("J132").Value = 3377.0715
$ws.Range("K132").Value = 95149284
$ws.Range("L132").Value = 10131.2145
$ws.Range("M132").Value = -95146754
$ws.Range("N132").Value = -15191.2145
$ws.Range("H135").Value = 1853668.6
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 1853668.6
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 16683017.4
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -16688087.4
$ws.Range("H137").Value = 596738.6
$ws.Range("I137").Value = 917004.9399999999
$ws.Range("J137").Value = 1958.4286
$ws.Range("K137").Value = 2751014.82
$ws.Range("L137").Value = 5875.2858
$ws.Range("M137").Value = -2748464.82
$ws.Range("N137").Value = -10975.2858
$ws.Range("H141").Value = 17433.125
$ws.Range("I141").Value = 25263.75
$ws.Range("J141").Value = 9602.5
$ws.Range("K141").Value = 75791.25
$ws.Range("L141").Value = 28807.5
$ws.Range("M141").Value = -70611.25
$ws.Range("N141").Value = -39167.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4646
$ws.Range("I31").Value = 988.8627300000001
$ws.Range("J31").Value = 25369.777
$ws.Range("K31").Value = 988.8627300000001
$ws.Range("L31").Value = 25369.777
$ws.Range("M31").Value = -693.8627300000001
$ws.Range("N31").Value = -25959.777
$ws.Range("H34").Value = 4646
$ws.Range("I34").Value = 988.8627300000001
$ws.Range("J34").Value = 25369.777
$ws.Range("K34").Value = 988.8627300000001
$ws.Range("L34").Value = 25369.777
$ws.Range("M34").Value = -786.8627300000001
$ws.Range("N34").Value = -25773.777
$ws.Range("H94").Value = 1864.9048
$ws.Range("I94").Value = 999.75
$ws.Range("J94").Value = 2068.4707
$ws.Range("K94").Value = 999.75
$ws.Range("L94").Value = 2068.4707
$ws.Range("M94").Value = -548.75
$ws.Range("N94").Value = -2970.4707
$ws.Range("H132").Value = 7411283.5
$ws.Range("I132").Value = 13889981
$ws.Range("J132").Value = 7057.6665
$ws.Range("K132").Value = 41669943
$ws.Range("L132").Value = 21172.9995
$ws.Range("M132").Value = -41667413
$ws.Range("N132").Value = -26232.9995
$ws.Range("H134").Value = 8447661
$ws.Range("I134").Value = 16668111
$ws.Range("J134").Value = 2842808.8
$ws.Range("K134").Value = 50004333
$ws.Range("L134").Value = 8528426.399999999
$ws.Range("M134").Value = -50001798
$ws.Range("N134").Value = -8533496.399999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 563.7059
$ws.Range("I5").Value = 434.81818
$ws.Range("J5").Value = 800
$ws.Range("K5").Value = 1304.45454
$ws.Range("L5").Value = 2400
$ws.Range("M5").Value = -1192.45454
$ws.Range("N5").Value = -2624
$ws.Range("H113").Value = 1540476
$ws.Range("I113").Value = 455.6757
$ws.Range("J113").Value = 3575502.8
$ws.Range("K113").Value = 1367.0271
$ws.Range("L113").Value = 10726508.4
$ws.Range("M113").Value = 802.9729
$ws.Range("N113").Value = -10730848.4
$ws.Range("H122").Value = 1462.0714
$ws.Range("I122").Value = 1474.5834
$ws.Range("J122").Value = 1387
$ws.Range("K122").Value = 13271.2506
$ws.Range("L122").Value = 12483
$ws.Range("M122").Value = -10821.2506
$ws.Range("N122").Value = -17383
$ws.Range("H126").Value = 1312.5
$ws.Range("I126").Value = 500
$ws.Range("K126").Value = 1500
$ws.Range("M126").Value = 3440
$ws.Range("H135").Value = 563.7059
$ws.Range("I135").Value = 434.81818
$ws.Range("J135").Value = 800
$ws.Range("K135").Value = 3913.36362
$ws.Range("L135").Value = 7200
$ws.Range("M135").Value = -1378.36362
$ws.Range("N135").Value = -12270

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1280.4736
$ws.Range("I102").Value = 1275.2667
$ws.Range("J102").Value = 1300
$ws.Range("K102").Value = 1275.2667
$ws.Range("L102").Value = 1300
$ws.Range("M102").Value = 346.7333000000001
$ws.Range("N102").Value = -4544

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1460.9565
$ws.Range("I7").Value = 1460.9565
$ws.Range("K7").Value = 1460.9565
$ws.Range("M7").Value = -1348.9565
$ws.Range("H16").Value = 6798.4346
$ws.Range("I16").Value = 1609.7368
$ws.Range("J16").Value = 31444.75
$ws.Range("K16").Value = 1609.7368
$ws.Range("L16").Value = 31444.75
$ws.Range("M16").Value = -1439.7368
$ws.Range("N16").Value = -31784.75
$ws.Range("H40").Value = 3634.1333
$ws.Range("I40").Value = 3690.7
$ws.Range("J40").Value = 3521
$ws.Range("K40").Value = 3690.7
$ws.Range("L40").Value = 3521
$ws.Range("M40").Value = -3554.7
$ws.Range("N40").Value = -3793
$ws.Range("H122").Value = 73542430
$ws.Range("I122").Value = 83351140
$ws.Range("J122").Value = 50001560
$ws.Range("K122").Value = 250053420
$ws.Range("L122").Value = 150004680
$ws.Range("M122").Value = -250050970
$ws.Range("N122").Value = -150009580
$ws.Range("H126").Value = 1460.9565
$ws.Range("I126").Value = 1460.9565
$ws.Range("K126").Value = 4382.8695
$ws.Range("M126").Value = -1912.8695
$ws.Range("H132").Value = 5408410
$ws.Range("I132").Value = 12502561
$ws.Range("J132").Value = 3343.1904
$ws.Range("K132").Value = 37507683
$ws.Range("L132").Value = 10029.5712
$ws.Range("M132").Value = -37505153
$ws.Range("N132").Value = -15089.5712
$ws.Range("H136").Value = 4627.7144
$ws.Range("I136").Value = 5563
$ws.Range("J136").Value = 3182.2727
$ws.Range("K136").Value = 16689
$ws.Range("L136").Value = 9546.8181
$ws.Range("M136").Value = -14139
$ws.Range("N136").Value = -14646.8181

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 19567.334
$ws.Range("I80").Value = 18700
$ws.Range("J80").Value = 20001
$ws.Range("K80").Value = 18700
$ws.Range("L80").Value = 20001
$ws.Range("M80").Value = -17702
$ws.Range("N80").Value = -21997
$ws.Range("H83").Value = 19567.334
$ws.Range("I83").Value = 18700
$ws.Range("J83").Value = 20001
$ws.Range("K83").Value = 56100
$ws.Range("L83").Value = 60003
$ws.Range("M83").Value = -51108
$ws.Range("N83").Value = -69987
$ws.Range("H126").Value = 83334200
$ws.Range("I126").Value = 27778260
$ws.Range("J126").Value = 166668110
$ws.Range("K126").Value = 83334780
$ws.Range("L126").Value = 500004330
$ws.Range("M126").Value = -83332310
$ws.Range("N126").Value = -500009270
$ws.Range("H132").Value = 17441096
$ws.Range("I132").Value = 8743195
$ws.Range("J132").Value = 28555082
$ws.Range("K132").Value = 26229585
$ws.Range("L132").Value = 85665246
$ws.Range("M132").Value = -26227055
$ws.Range("N132").Value = -85670306
$ws.Range("H136").Value = 15872912
$ws.Range("I136").Value = 8710183
$ws.Range("J136").Value = 38463060
$ws.Range("K136").Value = 26130549
$ws.Range("L136").Value = 115389180
$ws.Range("M136").Value = -26127999
$ws.Range("N136").Value = -115394280

Write-Output "Applied 192 cell changes"